# "Actualizar" automation run - appends one more availability-check pass
# (14 services) to the Disponibilidad log sheet, and corrects the
# floating-point precision of the previous pass's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up the timestamp precision on the previous run (rows 646-659) ---
$correctedPrevTimestamp = 44232.64054723379
for ($r = 646; $r -le 659; $r++) {
    $ws.Cells.Item($r, 4).Value = $correctedPrevTimestamp
}

# --- 2. Append the new run: rows 660-673 -------------------------------
# Each run is the same 14-service block cycled through; this mirrors the
# pattern already present at rows 2-15 / 646-659 etc.
$newRunTimestamp = 44232.66163265159

$services = @(
    [PSCustomObject]@{ Name = "Odoo";               Url = "https://www.dataintelligence-group.com/";                    SubAddress = $null },
    [PSCustomObject]@{ Name = "Blackbox";            Url = "https://serviciodashboard.azurewebsites.net/";               SubAddress = $null },
    [PSCustomObject]@{ Name = "PowerBI";             Url = "https://powerbi.microsoft.com/es-es/";                       SubAddress = $null },
    [PSCustomObject]@{ Name = "Dropbox";             Url = "https://www.dropbox.com/";                                   SubAddress = $null },
    [PSCustomObject]@{ Name = "Odoo";                Url = "https://dataintelligence.store/";                            SubAddress = $null },
    [PSCustomObject]@{ Name = "GEE";                 Url = "https://app-data-i.users.earthengine.app/";                  SubAddress = $null },
    [PSCustomObject]@{ Name = "UtilidadesOdoo";      Url = "https://odooutil.azurewebsites.net/";                        SubAddress = $null },
    [PSCustomObject]@{ Name = "Filtros Dashboard";   Url = "https://filtradordashboard.azurewebsites.net/";              SubAddress = $null },
    [PSCustomObject]@{ Name = "MapStore";            Url = "https://ide.dataintelligence-group.com/mapstore/";          SubAddress = "/" },
    [PSCustomObject]@{ Name = "GeoServer";           Url = "https://ide.dataintelligence-group.com/geoserver/web/?0";    SubAddress = $null },
    [PSCustomObject]@{ Name = "Tomcat";              Url = "https://ide.dataintelligence-group.com/";                    SubAddress = $null },
    [PSCustomObject]@{ Name = "Shiny";               Url = "https://rpubs.com/dataintelligence/";                       SubAddress = $null },
    [PSCustomObject]@{ Name = "Github";              Url = "https://github.com/Sud-Austral/";                           SubAddress = $null },
    [PSCustomObject]@{ Name = "EZ Exporter";         Url = "https://ezexporter.highviewapps.com/exports/export-profile/"; SubAddress = $null }
)

$startRow = 660
for ($i = 0; $i -lt $services.Count; $i++) {
    $r = $startRow + $i
    $svc = $services[$i]

    # The displayed URL text includes the sub-address fragment (e.g. "#/"),
    # same as every earlier block in this sheet.
    if ($svc.SubAddress) {
        $displayUrl = $svc.Url + "#" + $svc.SubAddress
    } else {
        $displayUrl = $svc.Url
    }

    $ws.Cells.Item($r, 1).Value = $svc.Name
    $ws.Cells.Item($r, 2).Value = $displayUrl
    $ws.Cells.Item($r, 3).Value = "Disponible"
    $ws.Cells.Item($r, 4).Value = $newRunTimestamp
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $linkCell = $ws.Cells.Item($r, 2)
    if ($svc.SubAddress) {
        $ws.Hyperlinks.Add($linkCell, $svc.Url, $svc.SubAddress)
    } else {
        $ws.Hyperlinks.Add($linkCell, $svc.Url)
    }
    $linkCell.Style = "Hyperlink"
}

Write-Output "Appended rows 660-673 and corrected D646:D659"
